# Apply the commit:
#  - Replace the "Fluentd" logging module shape with "Container Insights"
#  - Refresh the cached datetimeFigureOut footer field text
#    (6/4/19 -> 6/24/19) across the slide master, every slide layout,
#    and the notes master.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Swap the shape text: Fluentd -> Container Insights -----------------
$s.Shapes.Item("Rounded Rectangle 8").TextFrame.TextRange.Text = "Container Insights"

# --- 2. Refresh the cached date field everywhere it is rendered ------------
$newDate = "6/24/19"

# Slide master's own Date placeholder
$master = $s.Master
$master.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = $newDate

# Every slide layout owned by the master has its own Date placeholder copy
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $layout.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = $newDate
}

# Notes master keeps its own cached date field too
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate
